$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# The "survey" sheet's depth-to (column D) values were recalculated: every
# value in D2:D23 is shifted down by 90 (e.g. 40 -> -50, 30 -> -60, 42 -> -48).
for ($r = 2; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $old = $cell.Value2
    $cell.Value2 = $old - 90
}

# Update the sheet's active selection to match where the editor left off.
$ws.Range("B23").Select() | Out-Null
